$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 222, pushing the existing rows 222-226 down to 223-227
$ws.Rows.Item(222).Insert()

# Populate the newly inserted row 222 with the new record
$ws.Cells.Item(222, 1).Value = 3
$ws.Cells.Item(222, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(222, 3).Value = "Coquimbo"
$ws.Cells.Item(222, 4).Value = 44939
$ws.Cells.Item(222, 5).Value = 5
$ws.Cells.Item(222, 6).Value = 100112030
$ws.Cells.Item(222, 7).Value = "Poroto granado"
$ws.Cells.Item(222, 8).Value = "Sin especificar"
$ws.Cells.Item(222, 9).Value = "Primera"
$ws.Cells.Item(222, 10).Value = 73
$ws.Cells.Item(222, 11).Value = 39000
$ws.Cells.Item(222, 12).Value = 40000
$ws.Cells.Item(222, 13).Value = 39479
$ws.Cells.Item(222, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(222, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(222, 16).Value = 1579
$ws.Cells.Item(222, 17).Value = 25
$ws.Cells.Item(222, 18).Value = "Hortaliza"
